$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep only the header row and the first variant (demand1/net1/pv1/bat1)
# of each element/type pair (rows 1,2,4,6,8). Remove every other row
# (the "2" duplicate variants and the additional element types), i.e.
# delete rows 3,5,7,9-19, working from the bottom up so row numbers of
# rows still to be deleted don't shift.
$rowsToDelete = @(19,18,17,16,15,14,13,12,11,10,9,7,5,3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete() | Out-Null
}

$wb.Save()
